# The match-data for the two fixtures stored in rows 83 and 84 (ids 81 and 82)
# had been entered swapped; this corrects it by exchanging all the betting /
# result columns (B through AC) between the two rows while leaving column A
# (the "id" column) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($col in $cols) {
    $cell83 = $ws.Range($col + "83")
    $cell84 = $ws.Range($col + "84")

    $v83 = $cell83.Value2
    $v84 = $cell84.Value2

    if ($v83 -ne $v84) {
        $cell83.Value2 = $v84
        $cell84.Value2 = $v83
    }
}
